$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) counts per the source-site re-scrape.
# Mapping of sheet -> { row -> new value } built from the authoritative diff.

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 486  # was 485
$ws.Range("F4").Value = 25  # was 21
$ws.Range("F5").Value = 34  # was 33
$ws.Range("F6").Value = 53  # was 49
$ws.Range("F7").Value = 1129  # was 1125
$ws.Range("F9").Value = 224  # was 222
$ws.Range("F10").Value = 317  # was 310
$ws.Range("F11").Value = 7972  # was 7945
$ws.Range("F13").Value = 9443  # was 9392
$ws.Range("F14").Value = 73  # was 72
$ws.Range("F15").Value = 9  # was 8
$ws.Range("F17").Value = 463  # was 458
$ws.Range("F24").Value = 26  # was 25
$ws.Range("F25").Value = 41  # was 40
$ws.Range("F29").Value = 1620  # was 1613
$ws.Range("F30").Value = 23  # was 19
$ws.Range("F31").Value = 64  # was 61
$ws.Range("F32").Value = 299  # was 298
$ws.Range("F33").Value = 269  # was 268
$ws.Range("F34").Value = 41  # was 40
$ws.Range("F35").Value = 335  # was 332
$ws.Range("F36").Value = 55  # was 54
$ws.Range("F37").Value = 920  # was 913
$ws.Range("F41").Value = 402  # was 401
$ws.Range("F42").Value = 303  # was 302
$ws.Range("F43").Value = 265  # was 264
$ws.Range("F47").Value = 238  # was 236
$ws.Range("F48").Value = 87  # was 83

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 14  # was 12
$ws.Range("F15").Value = 47  # was 46
$ws.Range("F19").Value = 17  # was 15

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 202  # was 201
$ws.Range("F3").Value = 2729  # was 2726

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 202  # was 201
$ws.Range("F4").Value = 486  # was 485
$ws.Range("F8").Value = 26  # was 21
$ws.Range("F9").Value = 34  # was 33
$ws.Range("F10").Value = 53  # was 49
$ws.Range("F11").Value = 1129  # was 1125
$ws.Range("F14").Value = 14  # was 12
$ws.Range("F16").Value = 317  # was 310
$ws.Range("F17").Value = 7972  # was 7945
$ws.Range("F18").Value = 9443  # was 9392
$ws.Range("F19").Value = 73  # was 72
$ws.Range("F21").Value = 463  # was 458
$ws.Range("F24").Value = 26  # was 25
$ws.Range("F25").Value = 41  # was 40
$ws.Range("F28").Value = 1620  # was 1613
$ws.Range("F29").Value = 23  # was 19
$ws.Range("F30").Value = 64  # was 61
$ws.Range("F31").Value = 299  # was 298
$ws.Range("F32").Value = 269  # was 268
$ws.Range("F33").Value = 335  # was 332
$ws.Range("F34").Value = 55  # was 54
$ws.Range("F36").Value = 920  # was 913
$ws.Range("F39").Value = 402  # was 401
$ws.Range("F40").Value = 47  # was 46
$ws.Range("F41").Value = 303  # was 302
$ws.Range("F42").Value = 265  # was 264
$ws.Range("F46").Value = 238  # was 236
$ws.Range("F47").Value = 17  # was 15
$ws.Range("F49").Value = 87  # was 83
